# Commit before massive changes
# - Removes the former row 3 (T1234567J / Successful record) entirely.
# - Updates row 2 so it reflects a fresh registration snapshot:
#     C2: 2 -> 4, D2: "Successful" -> "Pending", E2: new timestamp.
# - Mirrors the row-header selection (rows 2:3) a user would make right
#   before deleting row 3, so the saved selection lands on A2:XFD3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select rows 2:3 (as a user would, to remove the now-duplicate record),
# then delete row 3 - this shifts nothing else and shrinks the used range
# to A1:E2.
$ws.Rows("2:3").Select()
$ws.Rows("3:3").Delete()

# Refresh the surviving row 2 with the latest registration data.
$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = "T2109876H"
$ws.Range("C2").Value = 4.0
$ws.Range("D2").Value = "Pending"
$ws.Range("E2").Value = 45768.5425244213
